$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL Scripts Assignment")

# --- New query-idea rows (5-16): Taken by / Progress / Description ---
# Column A ("Taken by") first, so "Nauman, O" becomes the first newly
# introduced shared string, then column B ("Progress"), then column C
# ("Description") - matching the order the strings were first used.
for ($r = 5; $r -le 16; $r++) {
    $ws.Range("A$r").Value = "Nauman, O"
}

$progress = @{
    5  = "Not Queried"
    6  = "Not Queried"
    7  = "Not Queried"
    8  = "Not Queried"
    9  = "Not Queried"
    10 = "Queried"
    11 = "Queried"
    12 = "Queried"
    13 = "In Progress"
    14 = "Not Queried"
    15 = "Queried"
    16 = "Queried"
}
for ($r = 5; $r -le 16; $r++) {
    $ws.Range("B$r").Value = $progress[$r]
}

$descriptions = @{
    5  = "Which countries have had the highest % of successful kickstarter projects?"
    6  = "Which Main Category has had the highest % of successful projects?"
    7  = "Which Sub Category within this has had the highest % of successful projects?"
    8  = "Which Kickstarter has raised the most `$?"
    9  = "Which Kickstarter has raised the most `$ per backer?"
    10 = "Avg kickstarter `$ raised per year"
    11 = "Avg kickstarter `$ raised per month"
    12 = "Avg goal"
    13 = "% of projects raising `$0 "
    14 = "Just an idea to play with: I think it would be cool to do something with the %Like% function in the title but am unsure of what would make the most sense? Maybe like something sad or dramatic?"
    15 = "Avg `$ raised per country"
    16 = "Number of campaigns per year"
}
for ($r = 5; $r -le 16; $r++) {
    $ws.Range("C$r").Value = $descriptions[$r]
}

# --- View-state updates: SQL Scripts Assignment becomes the active tab ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 60
$ws.Range("C18").Select() | Out-Null
